$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 167.44444
$ws.Range("J8").Value = 437.66666
$ws.Range("L8").Value = 1312.99998
$ws.Range("N8").Value = -1590.99998

$ws.Range("H12").Value = 292.5
$ws.Range("I12").Value = 297.22223
$ws.Range("K12").Value = 297.22223
$ws.Range("M12").Value = -127.22223

$ws.Range("H17").Value = 368.34286
$ws.Range("J17").Value = 311.87878
$ws.Range("L17").Value = 935.63634
$ws.Range("N17").Value = -1271.63634

$ws.Range("H31").Value = 20836816
$ws.Range("I31").Value = 22727436
$ws.Range("K31").Value = 68182308
$ws.Range("M31").Value = -68182078

$ws.Range("H32").Value = 12503626
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 12503626
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 12503626
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -12504278

$ws.Range("H69").Value = 25128.428
$ws.Range("I69").Value = 8250
$ws.Range("J69").Value = 31879.8
$ws.Range("K69").Value = 24750
$ws.Range("L69").Value = 95639.39999999999
$ws.Range("M69").Value = -23876
$ws.Range("N69").Value = -97387.39999999999

$ws.Range("H72").Value = 25128.428
$ws.Range("I72").Value = 8250
$ws.Range("J72").Value = 31879.8
$ws.Range("K72").Value = 74250
$ws.Range("L72").Value = 286918.2
$ws.Range("M72").Value = -69882
$ws.Range("N72").Value = -295654.2

$ws.Range("H87").Value = 29833.334
$ws.Range("J87").Value = 29833.334
$ws.Range("L87").Value = 29833.334
$ws.Range("N87").Value = -32329.334

$ws.Range("H90").Value = 29833.334
$ws.Range("J90").Value = 29833.334
$ws.Range("L90").Value = 89500.00199999999
$ws.Range("N90").Value = -101980.002

$ws.Range("H106").Value = 5756.4546
$ws.Range("I106").Value = 5744.579
$ws.Range("J106").Value = 5831.6665
$ws.Range("K106").Value = 5744.579
$ws.Range("L106").Value = 5831.6665
$ws.Range("M106").Value = -5113.579
$ws.Range("N106").Value = -7093.6665

$ws.Range("H113").Value = 8460.286
$ws.Range("I113").Value = 4751.3335
$ws.Range("J113").Value = 10217.158
$ws.Range("K113").Value = 4751.3335
$ws.Range("L113").Value = 10217.158
$ws.Range("M113").Value = -1497.3335
$ws.Range("N113").Value = -16725.158

$ws.Range("H116").Value = 7941
$ws.Range("I116").Value = 5059.4
$ws.Range("J116").Value = 9999.286
$ws.Range("K116").Value = 5059.4
$ws.Range("L116").Value = 9999.286
$ws.Range("M116").Value = -1617.4
$ws.Range("N116").Value = -16883.286

$ws.Range("H135").Value = 1254.4166
$ws.Range("I135").Value = 1217.1111
$ws.Range("K135").Value = 10953.9999
$ws.Range("M135").Value = -8418.999900000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 57598.25
$ws.Range("J28").Value = 200995
$ws.Range("L28").Value = 200995
$ws.Range("N28").Value = -201379

$ws.Range("H31").Value = 12251.917
$ws.Range("I31").Value = 11500.363
$ws.Range("K31").Value = 11500.363
$ws.Range("M31").Value = -11206.363

$ws.Range("H74").Value = 1132.9412
$ws.Range("I74").Value = 1132.9412
$ws.Range("K74").Value = 1132.9412
$ws.Range("M74").Value = -258.9412

$ws.Range("H77").Value = 1132.9412
$ws.Range("I77").Value = 1132.9412
$ws.Range("K77").Value = 5664.706
$ws.Range("M77").Value = -1296.706

$ws.Range("H99").Value = 57598.25
$ws.Range("J99").Value = 200995
$ws.Range("L99").Value = 200995
$ws.Range("N99").Value = -206985

$ws.Range("H105").Value = 53068
$ws.Range("J105").Value = 53068
$ws.Range("L105").Value = 53068
$ws.Range("N105").Value = -60056

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 660
$ws.Range("I22").Value = 500.33334
$ws.Range("J22").Value = 899.5
$ws.Range("K22").Value = 500.33334
$ws.Range("L22").Value = 899.5
$ws.Range("M22").Value = -327.33334
$ws.Range("N22").Value = -1245.5

$ws.Range("H102").Value = 33786.1
$ws.Range("I102").Value = 15207.333
$ws.Range("J102").Value = 200995
$ws.Range("K102").Value = 15207.333
$ws.Range("L102").Value = 200995
$ws.Range("M102").Value = -11962.333
$ws.Range("N102").Value = -207485

$ws.Range("H138").Value = 52352.53
$ws.Range("J138").Value = 52352.53
$ws.Range("L138").Value = 52352.53
$ws.Range("N138").Value = -62632.53

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1595.4166
$ws.Range("I58").Value = 1314.5
$ws.Range("K58").Value = 1314.5
$ws.Range("M58").Value = -1111.5

$ws.Range("H107").Value = 4832.12
$ws.Range("I107").Value = 671.3125
$ws.Range("J107").Value = 12229.111
$ws.Range("K107").Value = 671.3125
$ws.Range("L107").Value = 12229.111
$ws.Range("M107").Value = 1248.6875
$ws.Range("N107").Value = -16069.111

$ws.Range("H122").Value = 3436.5454
$ws.Range("J122").Value = 6006.25
$ws.Range("L122").Value = 18018.75
$ws.Range("N122").Value = -22918.75

$ws.Range("H132").Value = 5082
$ws.Range("I132").Value = 5082
$ws.Range("K132").Value = 15246
$ws.Range("M132").Value = -12716

$ws.Range("H134").Value = 2501.182
$ws.Range("I134").Value = 1808.5264
$ws.Range("K134").Value = 5425.5792
$ws.Range("M134").Value = -2890.5792

$ws.Range("H136").Value = 1595.4166
$ws.Range("I136").Value = 1314.5
$ws.Range("K136").Value = 3943.5
$ws.Range("M136").Value = -1393.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 83750856
$ws.Range("I4").Value = 91364120
$ws.Range("K4").Value = 274092360
$ws.Range("M4").Value = -274092248

$ws.Range("H12").Value = 562.9355
$ws.Range("J12").Value = 568.65216
$ws.Range("L12").Value = 1705.95648
$ws.Range("N12").Value = -2051.95648

$ws.Range("H21").Value = 175.85715
$ws.Range("I21").Value = 223.2
$ws.Range("J21").Value = 57.5
$ws.Range("K21").Value = 669.5999999999999
$ws.Range("L21").Value = 172.5
$ws.Range("M21").Value = -496.5999999999999
$ws.Range("N21").Value = -518.5

$ws.Range("H33").Value = 900
$ws.Range("I33").Value = 900
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 5400
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -5117
$ws.Range("N33").ClearContents()

$ws.Range("H46").Value = 100681.4
$ws.Range("I46").Value = 334066.34
$ws.Range("J46").Value = 659.2857
$ws.Range("K46").Value = 1002199.02
$ws.Range("L46").Value = 1977.8571
$ws.Range("M46").Value = -1002108.02
$ws.Range("N46").Value = -2159.8571

$ws.Range("H86").Value = 962.25
$ws.Range("I86").Value = 1082.25
$ws.Range("J86").Value = 722.25
$ws.Range("K86").Value = 3246.75
$ws.Range("L86").Value = 2166.75
$ws.Range("M86").Value = -2060.75
$ws.Range("N86").Value = -4538.75

$ws.Range("H87").Value = 1970
$ws.Range("I87").Value = 1970
$ws.Range("K87").Value = 5910
$ws.Range("M87").Value = -4662

$ws.Range("H89").Value = 962.25
$ws.Range("I89").Value = 1082.25
$ws.Range("J89").Value = 722.25
$ws.Range("K89").Value = 9740.25
$ws.Range("L89").Value = 6500.25
$ws.Range("M89").Value = -3812.25
$ws.Range("N89").Value = -18356.25

$ws.Range("H90").Value = 1970
$ws.Range("I90").Value = 1970
$ws.Range("K90").Value = 17730
$ws.Range("M90").Value = -11490

$ws.Range("H117").Value = 1184.5834
$ws.Range("I117").Value = 336
$ws.Range("K117").Value = 1008
$ws.Range("M117").Value = 2434

$ws.Range("H133").Value = 4109.3335
$ws.Range("I133").Value = 4109.3335
$ws.Range("K133").Value = 12328.0005
$ws.Range("M133").Value = -7268.000499999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 57228.24
$ws.Range("I80").Value = 141725.12
$ws.Range("K80").Value = 141725.12
$ws.Range("M80").Value = -140727.12

$ws.Range("H83").Value = 57228.24
$ws.Range("I83").Value = 141725.12
$ws.Range("K83").Value = 708625.6
$ws.Range("M83").Value = -703633.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 30000
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H98").Value = 29997.5
$ws.Range("J98").Value = 29997.5
$ws.Range("L98").Value = 29997.5
$ws.Range("N98").Value = -35987.5

$ws.Range("H123").Value = 39818.184
$ws.Range("J123").Value = 39818.184
$ws.Range("L123").Value = 39818.184
$ws.Range("N123").Value = -49618.184

$ws.Range("H136").Value = 2732.32
$ws.Range("I136").Value = 2097.9443
$ws.Range("J136").Value = 4363.5713
$ws.Range("K136").Value = 6293.8329
$ws.Range("L136").Value = 13090.7139
$ws.Range("M136").Value = -3743.8329
$ws.Range("N136").Value = -18190.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3535.8333
$ws.Range("J62").Value = 4291.2856
$ws.Range("L62").Value = 4291.2856
$ws.Range("N62").Value = -5539.2856

$ws.Range("H65").Value = 3535.8333
$ws.Range("J65").Value = 4291.2856
$ws.Range("L65").Value = 21456.428
$ws.Range("N65").Value = -27696.428
